# Swap the "Recorded By" entries from "dnasr281@gmail.com, System" to
# "System, dnasr281@gmail.com" wherever they occur in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldVal) {
        $cell.Value = $newVal
    }
}
